# Insert two new weekly price rows for "Chirimoya" (Vega Modelo de Temuco)
# right before the existing row 104, pushing all subsequent rows down by two
# (old row 104 becomes row 106, ..., old row 203 becomes row 205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 104 (each Insert() pushes the current
# row 104 and everything below it down by one).
$ws.Rows.Item(104).Insert()
$ws.Rows.Item(105).Insert()

# Seed the two new rows with a copy of the row that is now at 106 (the
# original row 104) so formatting/styles and the mostly-unchanged columns
# (A,B,C,E,F,G,H,I,J,K,Q,R,T) come along for free.
$ws.Rows.Item(106).Copy()
$ws.Rows.Item(104).PasteSpecial()
$ws.Rows.Item(106).Copy()
$ws.Rows.Item(105).PasteSpecial()
$excel.CutCopyMode = $false

# --- Row 104 -----------------------------------------------------------
$ws.Cells.Item(104, 4).Value = 45174       # D104 Fecha
$ws.Cells.Item(104, 12).Value = "Especial" # L104 Calidad
$ws.Cells.Item(104, 13).Value = 50         # M104 Volumen
$ws.Cells.Item(104, 14).Value = 3500       # N104 Precio minimo
$ws.Cells.Item(104, 15).Value = 3500       # O104 Precio maximo
$ws.Cells.Item(104, 16).Value = 3500       # P104 Precio promedio ponderado
$ws.Cells.Item(104, 19).Value = 3500       # S104 Precio $/Kg

# --- Row 105 -------------------------------------------------------------
$ws.Cells.Item(105, 4).Value = 45174       # D105 Fecha
$ws.Cells.Item(105, 13).Value = 120        # M105 Volumen
$ws.Cells.Item(105, 14).Value = 3000       # N105 Precio minimo
$ws.Cells.Item(105, 15).Value = 3000       # O105 Precio maximo
$ws.Cells.Item(105, 16).Value = 3000       # P105 Precio promedio ponderado
$ws.Cells.Item(105, 19).Value = 3000       # S105 Precio $/Kg
